$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FirstSet")

# New intake user credential row: username (hyperlinked mailto) + password.
$ws.Range("A7").Value = "intake.user1@cvhcare.com"
$ws.Range("B7").Value = "password"

# Add the mailto hyperlink for the new username cell (matches A2/A3/A4 pattern).
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:intake.user1@cvhcare.com")

# Re-apply the shared "Hyperlink" cell style so A7 matches A2/A3/A4 styling
# instead of the ad-hoc style variant Hyperlinks.Add creates.
$ws.Range("A7").Style = "Hyperlink"

# Reflect where editing left off.
$ws.Range("B7").Select() | Out-Null
